$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": the shared "In Translation" status text becomes
# "Handed back: in sync with en-US" for both the zh-cn and de-de columns.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.1666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.1666666666667

# ---------------------------------------------------------------------------
# Sheet "zh-cn": report the handed-back target/handback files for both rows.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/00b6ec3c0232a73d4efde654ec0237ff3ca6e0c9/e2e/12fccece-0abd-4a0a-8b82-deb48e7fc922.md", [Type]::Missing, [Type]::Missing, "12fccece-0abd-4a0a-8b82-deb48e7fc922.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/00b6ec3c0232a73d4efde654ec0237ff3ca6e0c9/e2e/12fccece-0abd-4a0a-8b82-deb48e7fc922.md", [Type]::Missing, [Type]::Missing, "12fccece-0abd-4a0a-8b82-deb48e7fc922.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/00b6ec3c0232a73d4efde654ec0237ff3ca6e0c9/e2e/315f1ec2-a1d7-4a34-863f-4c09b024cf56.md", [Type]::Missing, [Type]::Missing, "315f1ec2-a1d7-4a34-863f-4c09b024cf56.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/00b6ec3c0232a73d4efde654ec0237ff3ca6e0c9/e2e/315f1ec2-a1d7-4a34-863f-4c09b024cf56.md", [Type]::Missing, [Type]::Missing, "315f1ec2-a1d7-4a34-863f-4c09b024cf56.md")

$wsZh.Range("J2").Value = "12fccece-0abd-4a0a-8b82-deb48e7fc922.ce41db11c9657ea56624c5d79a07c3351891ecf7.zh-cn.xlf"
$wsZh.Range("J3").Value = "315f1ec2-a1d7-4a34-863f-4c09b024cf56.f9a08e67daf1d60948908658f31149f81861eeeb.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-13 18:29:16"
$wsZh.Range("K3").Value = "2016-08-13 18:29:16"

$wsZh.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsZh.Columns.Item(9).ColumnWidth = 39.1666666666667
$wsZh.Columns.Item(10).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------------
# Sheet "de-de": same handback reporting, with its own handback timestamp.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/00b6ec3c0232a73d4efde654ec0237ff3ca6e0c9/e2e/12fccece-0abd-4a0a-8b82-deb48e7fc922.md", [Type]::Missing, [Type]::Missing, "12fccece-0abd-4a0a-8b82-deb48e7fc922.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/00b6ec3c0232a73d4efde654ec0237ff3ca6e0c9/e2e/12fccece-0abd-4a0a-8b82-deb48e7fc922.md", [Type]::Missing, [Type]::Missing, "12fccece-0abd-4a0a-8b82-deb48e7fc922.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/00b6ec3c0232a73d4efde654ec0237ff3ca6e0c9/e2e/315f1ec2-a1d7-4a34-863f-4c09b024cf56.md", [Type]::Missing, [Type]::Missing, "315f1ec2-a1d7-4a34-863f-4c09b024cf56.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/00b6ec3c0232a73d4efde654ec0237ff3ca6e0c9/e2e/315f1ec2-a1d7-4a34-863f-4c09b024cf56.md", [Type]::Missing, [Type]::Missing, "315f1ec2-a1d7-4a34-863f-4c09b024cf56.md")

$wsDe.Range("J2").Value = "12fccece-0abd-4a0a-8b82-deb48e7fc922.ce41db11c9657ea56624c5d79a07c3351891ecf7.de-de.xlf"
$wsDe.Range("J3").Value = "315f1ec2-a1d7-4a34-863f-4c09b024cf56.f9a08e67daf1d60948908658f31149f81861eeeb.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-13 18:29:26"
$wsDe.Range("K3").Value = "2016-08-13 18:29:26"

$wsDe.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsDe.Columns.Item(9).ColumnWidth = 39.1666666666667
$wsDe.Columns.Item(10).ColumnWidth = 39.1666666666667
